# MMF Tracker 2026 - January sheet update
# Recomputed Interest Earned (D), Days Active (E) and Closing Balance (F)
# for rows that accrued an extra day of interest, plus the TOTAL row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("January")

# Row 5 - Berlin
$ws.Range("D5").Value = 17.588772879034934
$ws.Range("E5").Value = 25
$ws.Range("F5").Value = 6017.5887728790349

# Row 6 - Gandia
$ws.Range("D6").Value = 13.181320587279485
$ws.Range("E6").Value = 23
$ws.Range("F6").Value = 6013.1813205872795

# Row 7 - Denver
$ws.Range("D7").Value = 15.100680738591191
$ws.Range("E7").Value = 19
$ws.Range("F7").Value = 9015.1006807385911

# Row 9 - Oslo
$ws.Range("D9").Value = 0.96536075101544538
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 4000.9653607510154

# Row 12 - Lisbon
$ws.Range("D12").Value = 2.9066717855538466
$ws.Range("E12").Value = 12
$ws.Range("F12").Value = 4002.906671785554

# Row 14 - Paris
$ws.Range("D14").Value = 0.24638814609888296
$ws.Range("E14").Value = 5
$ws.Range("F14").Value = 510.24638814609887

# Row 15 - Marseille
$ws.Range("D15").Value = 1.690542979025897
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 2501.6905429790258

# Row 16 - Rio
$ws.Range("D16").Value = 14.628852760075757
$ws.Range("E16").Value = 18
$ws.Range("F16").Value = 6014.6288527600755

# Row 17 - Nairobi
$ws.Range("D17").Value = 19.302595666273302
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = 20019.302595666275

# Row 18 - Tel Aviv
$ws.Range("D18").Value = 29.274135303104856
$ws.Range("E18").Value = 24
$ws.Range("F18").Value = 14029.274135303105

# Row 19 - Valencia
$ws.Range("D19").Value = 72.719742301383704
$ws.Range("E19").Value = 11
$ws.Range("F19").Value = 60072.719742301386

# Row 20 - Manilla
$ws.Range("D20").Value = 11.732334614232885
$ws.Range("E20").Value = 22
$ws.Range("F20").Value = 4011.732334614233

# Row 21 - Scofield
$ws.Range("D21").Value = 1.2067009387693066
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 5001.2067009387692

# Row 23 - Emirates
$ws.Range("D23").Value = 13.147917729468245
$ws.Range("E23").Value = 24
$ws.Range("F23").Value = 7013.1479177294686

# Row 25 - Doha
$ws.Range("D25").Value = 4.8479828200922483
$ws.Range("E25").Value = 11
$ws.Range("F25").Value = 4004.8479828200921

# Row 26 - TOTAL
$ws.Range("D26").Value = 218.54
$ws.Range("F26").Value = 152228.54000000004

# Update the selected cell to match the saved view state
$ws.Range("H22").Select()
